# Tratamento de um novo formato do numero do processo
# Correcao na filtragem do despacho da patente
#
# Adds three new rows (39-41) to the "Itens de desenvolvimento" tracking
# sheet, mirroring the existing green "Defeito/Desenvolvido" rows (39-40)
# and the yellow "Melhoria/Desenvolvido" rows (41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row values ---------------------------------------------------
$ws.Range("A39").Value = "Filtragem incorreta ao buscar os despachos de patentes no cadastro do processo"
$ws.Range("B39").Value = "Defeito"
$ws.Range("C39").Value = "Desenvolvido"
$ws.Range("D39").Value = "N/A"

$ws.Range("A40").Value = "Filtragem incorreta das pastas cadastradas"
$ws.Range("B40").Value = "Defeito"
$ws.Range("C40").Value = "Desenvolvido"
$ws.Range("D40").Value = "N/A"

$ws.Range("A41").Value = "Implementado novo tratamento do número do processo na leitura da revista"
$ws.Range("B41").Value = "Melhoria"
$ws.Range("C41").Value = "Desenvolvido"
$ws.Range("D41").Value = "N/A"

# --- Formatting ---------------------------------------------------------
# Rows 39-40 use the same "green" fill + justified text used by the other
# Defeito/Desenvolvido rows above (e.g. rows 31-33).
$green = $ws.Range("A39:D40")
$green.Interior.Pattern = 1
$green.Interior.Color = 2359075
$green.Interior.PatternColor = 13421619
$green.HorizontalAlignment = -4130

# Row 41 uses the "yellow" fill used by the Melhoria/Desenvolvido rows
# above (e.g. rows 34-38): column C gets the slightly darker yellow
# variant, the rest get the lighter one.
$yellowMain = $ws.Range("A41")
$yellowMain.Interior.Pattern = 1
$yellowMain.Interior.Color = 65535
$yellowMain.Interior.PatternColor = 8421631
$yellowMain.HorizontalAlignment = -4130

$yellowB = $ws.Range("B41")
$yellowB.Interior.Pattern = 1
$yellowB.Interior.Color = 65535
$yellowB.Interior.PatternColor = 8421631
$yellowB.HorizontalAlignment = -4130

$yellowD = $ws.Range("D41")
$yellowD.Interior.Pattern = 1
$yellowD.Interior.Color = 65535
$yellowD.Interior.PatternColor = 8421631
$yellowD.HorizontalAlignment = -4130

$yellowC = $ws.Range("C41")
$yellowC.Interior.Pattern = 1
$yellowC.Interior.Color = 65535
$yellowC.Interior.PatternColor = 13421619
$yellowC.HorizontalAlignment = -4130

# --- View state -----------------------------------------------------
$ws.Range("A42").Select()
